# Update cryptos list with latest prices/volumes scraped on Tue Nov 28 14:49:11 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.477.05'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.60%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.036.96'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.82%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.24'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.13%  '

# Row 6
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.51%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.31'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.14%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.382'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.99%  '

# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.82%  '

# Row 11
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.23%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.341.28'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.92%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.44'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.31%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.35'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.55%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.745'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.22%  '

# Row 16
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.99%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.036.99'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.04%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.408.53'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.60%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.22'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.49%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.09'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.71%  '

# Row 21
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.86%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.84'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.37%  '

# Row 23
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.01%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.45'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.83%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.25'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.26%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.70'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.16%  '

# Row 27
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.04%  '

# Row 28
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +6.33%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.80'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.66%  '

# Row 30
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.31%  '

# Row 31
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.90%  '

# Row 32
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.80%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0607'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.13%  '

# Row 34
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.50'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.74%  '

# Row 35
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.02'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +10.08%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.34'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.06%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.76'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +9.09%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.22'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.44%  '

# Row 39
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.13%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.475.64'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.87%  '

# Row 41
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.29%  '

# Row 42
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.75%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '94.97'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.22%  '

# Row 44
$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.27'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +18.18%  '

# Row 45
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.81'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.39%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.32'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.36%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.12'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.28%  '

# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.93%  '

# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.60%  '

# Row 50
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.32%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.225.72'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.81%  '
